$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# 1) Insert a new worksheet "2022-Q4" right after "总计" (before "2022-Q2")
# ------------------------------------------------------------------
$q2SheetRef = $wb.Worksheets.Item("2022-Q2")
$newSheet = $wb.Worksheets.Add($q2SheetRef)
$newSheet.Name = "2022-Q4"

# Re-fetch sheet references by name (the pre-Add anchors go stale once the
# sheet collection is restructured by Add()).
$totalSheet = $wb.Worksheets.Item("总计")
$q2Sheet = $wb.Worksheets.Item("2022-Q2")

# --- Header row (row 1) ---
$newSheet.Range("B1").Value = "基金代码"
$newSheet.Range("C1").Value = "基金名称"
$newSheet.Range("D1").Value = "基金规模"
$newSheet.Range("E1").Value = "股票总仓位"
$newSheet.Range("F1").Value = "仓位占比"
$newSheet.Range("G1").Value = "持有市值(亿元)"
$newSheet.Range("H1").Value = "仓位排名"

$q2Sheet.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

# --- Data row (row 2) - fund holdings for 2022-Q4 ---
$newSheet.Range("A2").Value = 0
$newSheet.Range("B2").NumberFormat = "@"
$newSheet.Range("B2").Value = "165531"
$newSheet.Range("C2").Value = "信诚多策略灵活配置混合（LOF）"
$newSheet.Range("D2").NumberFormat = "@"
$newSheet.Range("D2").Value = "0.89"
$newSheet.Range("E2").NumberFormat = "@"
$newSheet.Range("E2").Value = "72.25"
$newSheet.Range("F2").NumberFormat = "@"
$newSheet.Range("F2").Value = "1.08"
$newSheet.Range("G2").NumberFormat = "@"
$newSheet.Range("G2").Value = "0.0096"
$newSheet.Range("H2").Value = 5

# Reset B2/D2/E2/F2/G2 back to the plain (unstyled) format used by the
# matching data cells on the reference sheet, now that the text values
# have been committed - mirrors source file (no explicit style on these).
$q2Sheet.Range("B2").Copy()
$newSheet.Range("B2").PasteSpecial(-4122)
$q2Sheet.Range("D2").Copy()
$newSheet.Range("D2").PasteSpecial(-4122)
$q2Sheet.Range("E2").Copy()
$newSheet.Range("E2").PasteSpecial(-4122)
$q2Sheet.Range("F2").Copy()
$newSheet.Range("F2").PasteSpecial(-4122)
$q2Sheet.Range("G2").Copy()
$newSheet.Range("G2").PasteSpecial(-4122)

$q2Sheet.Range("A2").Copy()
$newSheet.Range("A2").PasteSpecial(-4122)
$newSheet.Range("A2").Value = 0

# ------------------------------------------------------------------
# 2) Update "总计" sheet: shift rows 2-5 down to 3-6 and insert the new
#    "2022-Q4" summary row at row 2.
# ------------------------------------------------------------------
$totalSheet.Range("A6").Value = $totalSheet.Range("A5").Value()
$totalSheet.Range("B6").Value = $totalSheet.Range("B5").Value()
$totalSheet.Range("C6").Value = $totalSheet.Range("C5").Value()
$totalSheet.Range("D6").Value = $totalSheet.Range("D5").Value()

$totalSheet.Range("A5").Value = $totalSheet.Range("A4").Value()
$totalSheet.Range("B5").Value = $totalSheet.Range("B4").Value()
$totalSheet.Range("C5").Value = $totalSheet.Range("C4").Value()
$totalSheet.Range("D5").Value = $totalSheet.Range("D4").Value()

$totalSheet.Range("A4").Value = $totalSheet.Range("A3").Value()
$totalSheet.Range("B4").Value = $totalSheet.Range("B3").Value()
$totalSheet.Range("C4").Value = $totalSheet.Range("C3").Value()
$totalSheet.Range("D4").Value = $totalSheet.Range("D3").Value()

$totalSheet.Range("A3").Value = $totalSheet.Range("A2").Value()
$totalSheet.Range("B3").Value = $totalSheet.Range("B2").Value()
$totalSheet.Range("C3").Value = $totalSheet.Range("C2").Value()
$totalSheet.Range("D3").Value = $totalSheet.Range("D2").Value()

# New row 2 = 2022-Q4 summary
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.01

# Make sure the "index" column (A) keeps its original header-like style
# on the newly written/shifted rows (copy from A4, a row untouched in
# terms of style the whole time).
$totalSheet.Range("A4").Copy()
$totalSheet.Range("A2").PasteSpecial(-4122)
$totalSheet.Range("A6").PasteSpecial(-4122)

# Renumber the A column running index 0..4
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A5").Value = 3
$totalSheet.Range("A6").Value = 4

# ------------------------------------------------------------------
# 3) Restore the originally active tab ("2020-Q4"), which otherwise
#    loses its selection once a new sheet is inserted into the workbook.
# ------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item("2020-Q4")
$q4Sheet.Activate()
